$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Defined names -----------------------------------------------------
$wb.Names.Add("STR", "=TestSheet!`$B`$2")
$wb.Names.Add("DEX", "=TestSheet!`$C`$2")

# --- Header row strings (B1/C1) ----------------------------------------
$ws.Range("B1").Value = "STR"
$ws.Range("C1").Value = "DEX"

# --- Data values ---------------------------------------------------------
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 7

# --- Formula in A1 -------------------------------------------------------
$ws.Range("A1").Formula = "=SUM(STR, DEX)"

# --- Selection -------------------------------------------------------------
$ws.Range("E7").Select()
